$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 3 corresponds to 8ccf1249-3845-41e5-92d0-7ff41aef3685 file being handed off
# ("Ready for handoff" is the new status recorded after generating the handoff report)
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# Keep the datetime number format intact while refreshing the "Latest Handoff Date/Datetime" values
$wsOverview.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("D2").Value = "2016-03-20 17:40:24"
$wsOverview.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("D3").Value = "2016-03-20 17:40:24"

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("E2").Value = "2016-03-20 17:40:15"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("E3").Value = "2016-03-20 17:40:15"

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("E2").Value = "2016-03-20 17:40:24"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("E3").Value = "2016-03-20 17:40:24"
